$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the rows that were dropped from the table ---
# base_2t_separate (row 14), base_1_null (row 25), base_2_null (row 26)
# Delete bottom-up so earlier row numbers stay valid while deleting.
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(25).Delete()
$ws.Rows.Item(14).Delete()

# --- Update cells that changed value/content after the deletions shifted rows up ---

# base_2 (row 4): new text-similarity + new scores
$ws.Range("E4").Value = "Hadamard, Difference"
$ws.Range("L4").Value = 0.89
$ws.Range("M4").Value = 0.74
$ws.Range("N4").Value = 0.81

# base_allcomp (row 5): new text-similarity, LSTM/dropout info, dense layers, scores
$ws.Range("E5").Value = "Cosine, Hadamard, Difference"
$ws.Range("H5").Value = "25 / 0 / 0"
$ws.Range("I5").Value = "(32, 16)"
$ws.Range("J5").Value = 0.5
$ws.Range("L5").Value = 0.86
$ws.Range("M5").Value = 0.82
$ws.Range("N5").Value = 0.84

# base_2_allsim (row 8)
$ws.Range("E8").Value = "Cosine, Inverse_l1, Hadamard, Difference"
$ws.Range("H8").Value = "25 / 0 / 0"
$ws.Range("I8").Value = "(32, 16)"
$ws.Range("J8").Value = 0.5
$ws.Range("L8").Value = 0.81
$ws.Range("M8").Value = 0.74
$ws.Range("N8").Value = 0.78

# base_allcomp_allsim (row 9)
$ws.Range("E9").Value = "Cosine, Inverse_l1, Hadamard, Difference"
$ws.Range("H9").Value = "25 / 0 / 0"
$ws.Range("I9").Value = "(32, 16)"
$ws.Range("L9").Value = 0.87
$ws.Range("M9").Value = 0.73
$ws.Range("N9").Value = 0.79

# base_t_allcomp_allsim (row 10)
$ws.Range("E10").Value = "Cosine, Inverse_l1, Hadamard, Difference"
$ws.Range("H10").Value = "20 / 0.25 / 0.25"
$ws.Range("I10").Value = 24
$ws.Range("J10").Value = 0.75
$ws.Range("L10").Value = 0.83
$ws.Range("M10").Value = 0.77
$ws.Range("N10").Value = 0.8

# base_2t (row 13, after row14 deletion)
$ws.Range("E13").Value = "Hadamard, Difference"
$ws.Range("H13").Value = "25 / 0 / 0"
$ws.Range("I13").Value = "(32, 16)"
$ws.Range("J13").Value = 0.5
$ws.Range("L13").Value = 0.89
$ws.Range("M13").Value = 0.76
$ws.Range("N13").Value = 0.82

# base_2_full (row 16)
$ws.Range("E16").Value = "Hadamard, Difference"
$ws.Range("H16").Value = "25 / 0 / 0"
$ws.Range("I16").Value = "(32, 16)"
$ws.Range("J16").Value = 0.5
$ws.Range("L16").Value = 0.85
$ws.Range("M16").Value = 0.75
$ws.Range("N16").Value = 0.8

# base_2_num (row 23)
$ws.Range("E23").Value = "Hadamard, Difference"
$ws.Range("H23").Value = "25 / 0 / 0"
$ws.Range("I23").Value = "(32, 16)"
$ws.Range("J23").Value = 0.5
$ws.Range("L23").Value = 0.79
$ws.Range("M23").Value = 0.83
$ws.Range("N23").Value = 0.81

# base_2_num_null (row 25, after deletions)
$ws.Range("E25").Value = "Hadamard, Difference"
$ws.Range("H25").Value = "25 / 0 / 0"
$ws.Range("I25").Value = "(32, 16)"
$ws.Range("J25").Value = 0.5
$ws.Range("L25").Value = 0.85
$ws.Range("M25").Value = 0.83
$ws.Range("N25").Value = 0.84

# rows 27/28 swap order: reduced_all now above full_all
$ws.Range("A27").Value = "reduced_all"
$ws.Range("B27").Value = "Reduced"
$ws.Range("E27").Value = "Cosine, Inverse_l1, Hadamard, Difference"
$ws.Range("H27").Value = "25 / 0 / 0"
$ws.Range("I27").Value = 24
$ws.Range("J27").Value = 0.75
$ws.Range("L27").Value = 0.87
$ws.Range("M27").Value = 0.82
$ws.Range("N27").Value = 0.85

$ws.Range("A28").Value = "full_all"
$ws.Range("B28").Value = "Full"
$ws.Range("E28").Value = "Cosine, Inverse_l1, Hadamard, Difference"
$ws.Range("L28").Value = 0.84
$ws.Range("M28").Value = 0.8
$ws.Range("N28").Value = 0.82

# Row 27 picked up an explicit row height when its content moved
$ws.Rows.Item(27).RowHeight = 17

# --- Sheet-level view / layout tweaks ---
$ws.Columns.Item(5).ColumnWidth = 37.5

$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Application.ActiveWindow.Zoom = 82
$ws.Range("M36").Select()

# --- Workbook window geometry ---
$excel.ActiveWindow.Left = 52800
$excel.ActiveWindow.Top = -12860
$excel.ActiveWindow.Width = 19200
$excel.ActiveWindow.Height = 10800
